$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly driver report update: the two "Bad Drivers" rows swap order and get
# refreshed counts, and a missing "Driver Vintage" date is filled in.

# Row 3 now holds the 23.90.0.2 driver with its refreshed stats
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 363
$ws.Range("D3").Value = 98.7

# Row 4 now holds the 23.10.0.8 driver with its refreshed stats
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 98.8

# Totals row: Critical Minutes total refreshed (Client Count total unchanged)
$ws.Range("C5").Value = 369

# Fill in the previously-blank Driver Vintage date for row 13. Assign it
# through an apostrophe-prefixed literal so it is stored as text (like its
# siblings in E14/E15) rather than being auto-converted to a date serial,
# then copy E14's number format over so the cell keeps its original style.
$ws.Range("E13").Value = "'2022-08-29"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
